$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for the Price/Volume columns so that numeric-looking
# strings (e.g. "21.00", "1.176") are preserved exactly as text instead of
# being auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "23.830.33"
$ws.Range("E2").Value = "  +15.83%  "
$ws.Range("D3").Value = "1.658.58"
$ws.Range("E3").Value = "  +12.83%  "
$ws.Range("D4").Value = "0.9879"
$ws.Range("E4").Value = "  -2.15%  "
$ws.Range("D5").Value = "305.86"
$ws.Range("E5").Value = "  +10.41%  "
$ws.Range("D6").Value = "0.9796"
$ws.Range("E6").Value = "  +2.09%  "
$ws.Range("D7").Value = "0.3736"
$ws.Range("E7").Value = "  +4.92%  "
$ws.Range("D8").Value = "0.3448"
$ws.Range("E8").Value = "  +12.68%  "
$ws.Range("B9").Value = "OKB"
$ws.Range("C9").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D9").Value = "43.49"
$ws.Range("E9").Value = "  +10.31%  "
$ws.Range("B10").Value = "Polygon"
$ws.Range("C10").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D10").Value = "1.176"
$ws.Range("E10").Value = "  +8.51%  "
$ws.Range("D11").Value = "0.07223"
$ws.Range("E11").Value = "  +9.07%  "
$ws.Range("D12").Value = "0.9840"
$ws.Range("E12").Value = "  -1.96%  "
$ws.Range("D13").Value = "21.00"
$ws.Range("E13").Value = "  +16.18%  "
$ws.Range("D14").Value = "6.002"
$ws.Range("E14").Value = "  +10.06%  "
$ws.Range("D15").Value = "6.751"
$ws.Range("E15").Value = "  +9.49%  "
$ws.Range("D16").Value = "1.660.41"
$ws.Range("E16").Value = "  +12.87%  "
$ws.Range("D17").Value = "0.00001101"
$ws.Range("E17").Value = "  +8.14%  "
$ws.Range("D18").Value = "0.9790"
$ws.Range("E18").Value = "  +2.00%  "
$ws.Range("D19").Value = "0.06717"
$ws.Range("E19").Value = "  +12.70%  "
$ws.Range("D20").Value = "80.83"
$ws.Range("E20").Value = "  +17.12%  "
$ws.Range("D21").Value = "16.46"
$ws.Range("E21").Value = "  +13.99%  "
$ws.Range("D22").Value = "6.086"
$ws.Range("E22").Value = "  +11.12%  "
$ws.Range("D23").Value = "11.99"
$ws.Range("E23").Value = "  +6.92%  "
$ws.Range("D24").Value = "23.839.11"
$ws.Range("E24").Value = "  +15.86%  "
$ws.Range("D25").Value = "2.346"
$ws.Range("E25").Value = "  +2.92%  "
$ws.Range("B26").Value = "LEO"
$ws.Range("C26").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D26").Value = "3.404"
$ws.Range("E26").Value = "  -8.41%  "
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").Value = "2.702"
$ws.Range("E27").Value = "  +29.65%  "
$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").Value = "152.26"
$ws.Range("E28").Value = "  +4.82%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "19.55"
$ws.Range("E29").Value = "  +14.44%  "
$ws.Range("B30").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C30").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D30").Value = "1.837.96"
$ws.Range("E30").Value = "  +12.68%  "
$ws.Range("B31").Value = "BitcoinCash"
$ws.Range("C31").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D31").Value = "126.28"
$ws.Range("E31").Value = "  +10.22%  "
$ws.Range("D32").Value = "4.085"
$ws.Range("E32").Value = "  +3.63%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "6.191"
$ws.Range("E33").Value = "  +26.17%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "1.001"
$ws.Range("E34").Value = "  +26.51%  "
$ws.Range("B35").Value = "WEMIXTOKEN"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").Value = "1.702"
$ws.Range("E35").Value = "  +18.04%  "
$ws.Range("B36").Value = "Stellar"
$ws.Range("C36").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D36").Value = "0.08361"
$ws.Range("E36").Value = "  +5.22%  "
$ws.Range("D37").Value = "12.32"
$ws.Range("E37").Value = "  +20.21%  "
$ws.Range("B38").Value = "FraxShare"
$ws.Range("C38").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D38").Value = "8.920"
$ws.Range("E38").Value = "  +22.64%  "
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").Value = "5.303"
$ws.Range("E39").Value = "  +12.79%  "
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").Value = "0.06293"
$ws.Range("E40").Value = "  +10.44%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "1.277"
$ws.Range("E41").Value = "  +4.29%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "0.02283"
$ws.Range("E42").Value = "  +12.58%  "
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").Value = "0.2054"
$ws.Range("E43").Value = "  +11.27%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "0.6062"
$ws.Range("E44").Value = "  +16.32%  "
$ws.Range("B45").Value = "Frax"
$ws.Range("C45").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D45").Value = "0.9782"
$ws.Range("E45").Value = "  +1.85%  "
$ws.Range("D46").Value = "13.31"
$ws.Range("E46").Value = "  +11.17%  "
$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").Value = "3.818"
$ws.Range("E47").Value = "  +8.66%  "
$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").Value = "0.5936"
$ws.Range("E48").Value = "  +15.06%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "127.08"
$ws.Range("E49").Value = "  +6.05%  "
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "2.001"
$ws.Range("E50").Value = "  +11.43%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "0.07092"
$ws.Range("E51").Value = "  +10.31%  "
